$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell E8 text from "Good Morning" to "GIT UPDATE"
$ws.Range("E8").Value = "GIT UPDATE"

# Select E8 as the active cell (reflected in sheetView selection)
$ws.Range("E8").Select()
